$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "Ativação:" date, 01/01/2018 -> 01/01/2021 ---
# A direct Range.Value assignment of a date-shaped string ("01/01/2021") gets
# auto-recognized as an actual date by Excel's type inference and would turn
# the cell into a numeric date serial (changing both its stored type and its
# cell style/number format). To keep it as literal text (matching the
# original authoring, which stored the date as plain text), stage the text
# in a scratch cell via a formula (so it round-trips as a string result),
# then Copy / PasteSpecial the values into the target cells. Paste of a
# string value does not re-run the "looks like a date" autocorrection.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="01/01/2021"'
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()

# --- Row 13: "Docentes responsáveis:" ---
$ws.Range("B13").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C13").Value = "198273 - Domingos Savio Giordani"

# --- Row 19: "Método:" ---
$metodo = "Desenvolvimento e apresentação do Projeto monografia a ser desenvolvida na disciplina de Trabalho de Graduação em Engenharia de Produção II, conforme norma do Curso de Engenharia de Produção"
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Row 20: "Critério:" ---
$criterio = "Avaliação Ad hoc por 2 examinadores. A nota da disciplina será a média das duas notas"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Row 21: "Norma de recuperação:" ---
$norma = "Reapresentação do trabalho modificado para nova avaliação"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- Row 22: "Bibliografia:" ---
$bib = "Cauchick-Miguel, P. A.,   Metodologia de pesquisa em engenharia de produção e gestão de operações / Afonso Fleury ... [et al.] ; coordenação . - 3. ed. - Rio de Janeiro : Elsevier, 2018. Cauchick-Miguel, P. A.,   Metodologia de pesquisa em engenharia , 1. ed. - Rio de Janeiro : GEN LTC, 2019. BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005.GIL, A.C. Como elaborar projetos de pesquisa. 5 ed. Atlas, São Paulo, 2010.MEDEIROS, J. B. Redação Cientifica: A Prática de Fichamentos, Resumos e Resenhas. 11 ed. São Paulo: Atlas, 2009"
$ws.Range("B22").Value = $bib
$ws.Range("C22").Value = $bib
